# daily auto push: 2026-02-25 10:06 UTC
# Insert a new timetable entry for 2026/02/25 (水, 16時) before row 872,
# shifting the remaining rows (old 872:913) down to (873:914).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 872; everything below shifts down by one.
$ws.Rows("872:872").Insert()

# Fill in the new row's values. Column A holds dates stored as plain text
# (matching the rest of the sheet) rather than Excel's auto-detected date
# serials, so we briefly force a text format, assign the value, then clear
# the temporary formatting back off the cell.
$ws.Range("A872").NumberFormat = "@"
$ws.Range("A872").Value = "2026/02/25"
$ws.Range("A872").ClearFormats()

$ws.Range("B872").Value = "水"
$ws.Range("C872").Value = 16
$ws.Range("D872").Value = 201
